$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# ---- Header row (row 1) ----
$headers = @(
    "Date","Model Name","Exact Precision (Micro Avg)","Exact Recall (Micro Avg)",
    "Exact F1 Score (Micro Avg)","Exact Precision (Macro Avg)","Exact Recall (Macro Avg)",
    "Exact F1 Score (Macro Avg)","Exact Precision (Weighted Avg)","Exact Recall (Weighted Avg)",
    "Exact F1 Score (Weighted Avg)","Partial Precision","Partial Recall","Partial F1 Score",
    "Partial TP","Partial FP","Partial FN","Support","Accuracy","Result Link","Stats Link",
    "No of GPU Used","Power Consumption"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# Style the header row: bold font, thin box border, centered horizontal, top vertical alignment
$headerRange = $ws.Range("A1:W1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- Data row (row 2) ----
$ws.Cells.Item(2, 2).Value = "Llama-3.3-70B-Instruct"
$ws.Cells.Item(2, 3).Value = 0.3541666666666667
$ws.Cells.Item(2, 4).Value = 0.2687747035573123
$ws.Cells.Item(2, 5).Value = 0.3056179775280899
$ws.Cells.Item(2, 6).Value = 0.1646301846301846
$ws.Cells.Item(2, 7).Value = 0.1003718962324802
$ws.Cells.Item(2, 8).Value = 0.1239091829270396
$ws.Cells.Item(2, 9).Value = 0.426223374049461
$ws.Cells.Item(2, 10).Value = 0.2687747035573123
$ws.Cells.Item(2, 11).Value = 0.3268060112764173
$ws.Cells.Item(2, 12).Value = 0.4210526315789473
$ws.Cells.Item(2, 13).Value = 0.316205533596838
$ws.Cells.Item(2, 14).Value = 0.3611738148984199
$ws.Cells.Item(2, 15).Value = 80
$ws.Cells.Item(2, 16).Value = 110
$ws.Cells.Item(2, 17).Value = 173
$ws.Cells.Item(2, 18).Value = 253
$ws.Cells.Item(2, 19).Value = 0.9437826844857641
$ws.Cells.Item(2, 20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.3-70B-Instruct_5_shot.txt"
$ws.Cells.Item(2, 21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.3-70B-Instruct_5_shot.txt"
$ws.Cells.Item(2, 22).Value = "4 MLGPU"
$ws.Cells.Item(2, 23).Value = "0.142 kWh"
$ws.Cells.Item(2, 24).Value = 5647

# Date column needs special handling to avoid Excel auto-converting the
# date-like text "09/10/2025" into a date serial number: force text
# format first, then reset back to Normal so no stray style sticks to
# the cell.
$dateCell = $ws.Cells.Item(2, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/10/2025"
$dateCell.Style = "Normal"
